$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "triggers" worksheet between "static_Tables" and
#    "dynamic_table".
# ---------------------------------------------------------------------------
$staticTables = $wb.Worksheets.Item("static_Tables")
$triggers = $wb.Worksheets.Add($null, $staticTables)
$triggers.Name = "triggers"

# Column widths (characters)
$triggers.Columns.Item(1).ColumnWidth = 33.666666666666664
$triggers.Columns.Item(2).ColumnWidth = 23.333333333333336
$triggers.Columns.Item(3).ColumnWidth = 42.666666666666664

# Header row
$triggers.Range("A1:C1").Font.Bold = $true
$triggers.Range("A1").Value = "triggers"
$triggers.Range("B1").Value = "tables"
$triggers.Range("C1").Value = "Function"

# Data rows
$triggers.Range("A3").Value = "trigger_update_amt"
$triggers.Range("B3").Value = "order_item_table"
$triggers.Range("C3").Value = "This will calculate the amount of item wrt price of product and quantity selected."

$triggers.Range("A5").Value = "trigger_update_total_order_amount"
$triggers.Range("B5").Value = "order_item_table"
$triggers.Range("C5").Value = "Updates total_amount in order_table based onsum all items for that has same order_id"

$triggers.Range("A7").Value = "trigger_update_order_history"
$triggers.Range("B7").Value = "order_table"
$triggers.Range("C7").Value = "If there is an update in order_status it's been recorded in order_history_table"

$triggers.Range("A9").Value = "check_duplicate_registration"
$triggers.Range("B9").Value = "business_table"
$triggers.Range("C9").Value = "Checks if the b_registrationnumber is unique before inserting new "

$triggers.Range("A11").Value = "trigger_insert_business_user"
$triggers.Range("B11").Value = "business_table"
$triggers.Range("C11").Value = "checks if active_stautus = 1 in business_table , then enter the username , password in business_user_table"

$triggers.Range("A14").Value = "trigger_insert_user"
$triggers.Range("B14").Value = "user_table"
$triggers.Range("C14").Value = "checks if active_stautus = 1 in user_table , then enter the username , password in users."

$triggers.Range("B18").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. static_Tables sheet tweaks.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("static_Tables")

# B74 loses its style (goes back to the default/no style).
$ws1.Range("B74").ClearFormats() | Out-Null

# Remove the leftover blank/styled rows 255-267 and 269, keep row 268 (drop
# its B cell only).
$ws1.Range("A255:B267").Clear() | Out-Null
$ws1.Range("B268").Clear() | Out-Null
$ws1.Range("A269:B269").Clear() | Out-Null

$ws1.Range("B258").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. roles_access_to_each_table sheet: scroll position tweak.
# ---------------------------------------------------------------------------
$roles = $wb.Worksheets.Item("roles_access_to_each_table")
$roles.Activate()
$excel.ActiveWindow.ScrollRow = 32
$roles.Range("C70").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Leave "triggers" as the active sheet/tab when saved.
# ---------------------------------------------------------------------------
$triggers.Activate()
$triggers.Range("B18").Select() | Out-Null
